$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: quality_comparison
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("quality_comparison")

# C1 / D1 sit inside the merged header cell B1:D1. They keep their (empty)
# numeric value but lose the bold/boxed "Normal-table" look in favour of a
# plain top+bottom (C1) / top+bottom+right (D1) thin rule - matching the
# border used to underline the merged header band.
$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$c1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

$d1 = $ws1.Range("D1")
$d1.Style = "Normal"
$d1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$d1.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$d1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

# Anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"

# ---------------------------------------------------------------------
# Sheet 2: computational_comparison
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("computational_comparison")

$c1b = $ws2.Range("C1")
$c1b.Style = "Normal"
$c1b.Borders.Item(8).LineStyle = 1
$c1b.Borders.Item(9).LineStyle = 1

$d1b = $ws2.Range("D1")
$d1b.Style = "Normal"
$d1b.Borders.Item(8).LineStyle = 1
$d1b.Borders.Item(10).LineStyle = 1
$d1b.Borders.Item(9).LineStyle = 1

$f1b = $ws2.Range("F1")
$f1b.Style = "Normal"
$f1b.Borders.Item(8).LineStyle = 1
$f1b.Borders.Item(9).LineStyle = 1

$g1b = $ws2.Range("G1")
$g1b.Style = "Normal"
$g1b.Borders.Item(8).LineStyle = 1
$g1b.Borders.Item(10).LineStyle = 1
$g1b.Borders.Item(9).LineStyle = 1

# Anonymize "fedcore" -> "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Drop the stray empty inline-string cell at G5
$ws2.Range("G5").ClearContents()
